# Weekly fruit/vegetable price update: insert 3 new rows of "Palta" (avocado)
# price data for Comercializadora del Agro de Limari, ahead of the existing
# history, which shifts down by 3 rows (old row 422 -> 425, ... old row 442 -> 445).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows right before row 422, pushing rows 422:442 down to 425:445.
$ws.Rows("422:424").Insert()

# New week's data (date serial 44714 = 2022-06-02) for the three quality grades
# that were previously reported at rows 422:424 (now pushed to 425:427).
$newRows = @(
    @{ Row = 422; K = "Hass"; L = "Especial"; M = 240; N = 3200; O = 3300; P = 3250; S = 3250 },
    @{ Row = 423; K = "Hass"; L = "Primera";  M = 300; N = 3000; O = 3100; P = 3050; S = 3050 },
    @{ Row = 424; K = "Hass"; L = "Segunda";  M = 240; N = 2800; O = 2900; P = 2850; S = 2850 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = 2
    $ws.Cells.Item($row, 2).Value = "Comercializadora del Agro de Limarí"
    $ws.Cells.Item($row, 3).Value = "Coquimbo"

    $dcell = $ws.Cells.Item($row, 4)
    $dcell.Value = 44714
    $dcell.NumberFormat = $ws.Cells.Item(426, 4).NumberFormat

    $ws.Cells.Item($row, 5).Value = 4
    $ws.Cells.Item($row, 6).Value = "Fruta"
    $ws.Cells.Item($row, 7).Value = 100106
    $ws.Cells.Item($row, 8).Value = "Oleaginosos"
    $ws.Cells.Item($row, 9).Value = 100106002
    $ws.Cells.Item($row, 10).Value = "Palta"
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = "`$/kilo (en caja de 17 kilos)"
    $ws.Cells.Item($row, 18).Value = "Provincia de Limarí"
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = 1
}
